$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds crypto prices stored as literal text in the source data
# (values like "65.617.83" use dots as thousands separators and must not
# be reinterpreted as numbers), so force text format before writing them.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.617.83'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.648.37'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.98'
$ws.Range('E5').Value = '  -0.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.33'
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.628'
$ws.Range('E8').Value = '  +1.41%  '
$ws.Range('E9').Value = '  +2.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.82'
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.396'
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('E14').Value = '  -3.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.123.82'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.426.30'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.648.60'
$ws.Range('E17').Value = '  -0.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.57'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.73'
$ws.Range('E19').Value = '  -2.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.42'
$ws.Range('E20').Value = '  -2.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '348.36'
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.88'
$ws.Range('E23').Value = '  -1.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000113'
$ws.Range('E24').Value = '  +1.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.60'
$ws.Range('E25').Value = '  -2.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.69'
$ws.Range('E26').Value = '  +3.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.59'
$ws.Range('E27').Value = '  -1.44%  '
$ws.Range('E28').Value = '  -2.87%  '
$ws.Range('E29').Value = '  +0.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.89'
$ws.Range('E30').Value = '  -3.14%  '
$ws.Range('E31').Value = '  -1.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '526.60'
$ws.Range('E32').Value = '  -2.82%  '
$ws.Range('E33').Value = '  -1.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.39'
$ws.Range('E34').Value = '  -3.19%  '
$ws.Range('E35').Value = '  -1.76%  '
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.34'
$ws.Range('E37').Value = '  -0.52%  '
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.93'
$ws.Range('E39').Value = '  -1.48%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '155.73'
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '161.00'
$ws.Range('E42').Value = '  -2.75%  '
$ws.Range('E43').Value = '  -0.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0605'
$ws.Range('E44').Value = '  -1.66%  '
$ws.Range('E45').Value = '  -1.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.52'
$ws.Range('E46').Value = '  -3.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.635'
$ws.Range('E47').Value = '  -1.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0254'
$ws.Range('E48').Value = '  -2.39%  '
$ws.Range('E49').Value = '  -2.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₆0251'
$ws.Range('E50').Value = '  +6.74%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.70'
$ws.Range('E51').Value = '  -2.36%  '
